$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(22).EntireRow.Delete()
$ws.Rows.Item(22).EntireRow.Delete()
